# Auto-generated edit script: apply cell value changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.611.68"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "'2.244.60"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'306.42"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'94.48"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").Value = "'34.79"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "'0.0800"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "'7.18"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'2.587.94"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "'2.242.69"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").Value = "'0.830"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'13.55"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'44.397.18"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "'0.0₃0933"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "'11.74"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").Value = "'65.29"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'237.49"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'2.94"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +4.06%  "
$ws.Range("D28").Value = "'9.77"
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("D29").Value = "'36.92"
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("D30").Value = "'19.95"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "'5.84"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'147.96"
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "'3.18"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("E38").Value = "  +5.03%  "
$ws.Range("D39").Value = "'15.14"
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "'1.808.69"
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("E45").Value = "  +12.03%  "
$ws.Range("D46").Value = "'81.93"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "'98.11"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'4.81"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'68.77"
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("D51").Value = "'53.94"
$ws.Range("E51").Value = "  -1.02%  "
